$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "problèmes liés"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "problème lié"
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
